# "bug with empty notes is resolved"
#
# Two fixes to the footnotes of this document:
#   1. Footnote 24's text had a stray "aa" typo appended to the end of
#      "...པེ་ཅིན།" - strip it off.
#   2. Footnote 29 was an "empty" note (its body was just the punctuation
#      mark "།", with no real content) left over from editing; remove the
#      footnote entirely (both its reference in the body text and its
#      definition), which is the "bug with empty notes" the commit fixes.

$d = $word.ActiveDocument

# --- 1. Fix the "aa" typo trailing footnote 24's text -----------------
$footnotes = $d.Footnotes
for ($i = 1; $i -le $footnotes.Count; $i++) {
    $fn = $footnotes.Item($i)
    $text = $fn.Range.Text
    if ($text -like "*aa") {
        $fn.Range.Text = $text.Substring(0, $text.Length - 2)
    }
}

# --- 2. Remove the empty footnote (body text is just "།") -------------
for ($i = $footnotes.Count; $i -ge 1; $i--) {
    $fn = $footnotes.Item($i)
    $text = $fn.Range.Text
    if ($text -eq "།") {
        $fn.Delete()
    }
}

Write-Output "Footnotes remaining: $($d.Footnotes.Count)"
